$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source data stores every Price cell (column D) as text (inline string), even when the
# text looks like a plain number. Force the NumberFormat to Text ("@") on such cells right
# before assigning their Value so Excel does not auto-convert them into numeric cells.

$ws.Range("D2").Value = "20.281.16"
$ws.Range("E2").Value = "  +1.38%  "

$ws.Range("D3").Value = "1.450.57"
$ws.Range("E3").Value = "  +2.78%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.62%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9375"
$ws.Range("E5").Value = "  -6.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "273.59"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3638"
$ws.Range("E7").Value = "  -0.89%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3058"
$ws.Range("E8").Value = "  -1.96%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.72"
$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("E10").Value = "  -0.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06519"
$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9981"
$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.356"
$ws.Range("E13").Value = "  -2.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.64"
$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.074"
$ws.Range("E15").Value = "  -1.74%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001014"
$ws.Range("E16").Value = "  -0.62%  "

$ws.Range("D17").Value = "1.446.44"
$ws.Range("E17").Value = "  +2.37%  "

$ws.Range("E18").Value = "  -4.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05703"
$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.93"
$ws.Range("E20").Value = "  -2.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.392"
$ws.Range("E21").Value = "  -4.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.31"
$ws.Range("E22").Value = "  -2.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.80"
$ws.Range("E23").Value = "  -2.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.235"
$ws.Range("E24").Value = "  -1.13%  "

$ws.Range("D25").Value = "20.307.23"
$ws.Range("E25").Value = "  +1.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.74"
$ws.Range("E26").Value = "  +5.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.085"
$ws.Range("E27").Value = "  -7.90%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.98"
$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("D29").Value = "1.599.47"
$ws.Range("E29").Value = "  +1.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "111.19"
$ws.Range("E30").Value = "  +1.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.957"
$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.813"
$ws.Range("E32").Value = "  -9.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7840"
$ws.Range("E33").Value = "  -4.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07729"
$ws.Range("E34").Value = "  +0.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.486"
$ws.Range("E35").Value = "  +0.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05622"
$ws.Range("E36").Value = "  -4.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.648"
$ws.Range("E37").Value = "  -5.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.117"
$ws.Range("E38").Value = "  +2.34%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02004"
$ws.Range("E39").Value = "  -3.40%  "

$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9469"
$ws.Range("E40").Value = "  -5.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.19"
$ws.Range("E41").Value = "  -3.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1849"
$ws.Range("E42").Value = "  -2.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.355"
$ws.Range("E43").Value = "  -12.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5228"
$ws.Range("E44").Value = "  -1.59%  "

$ws.Range("E45").Value = "  -1.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.88"
$ws.Range("E46").Value = "  -3.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "116.58"
$ws.Range("E47").Value = "  +0.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5106"
$ws.Range("E48").Value = "  -1.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.735"
$ws.Range("E49").Value = "  -1.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06381"
$ws.Range("E50").Value = "  +2.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9843"
$ws.Range("E51").Value = "  -1.66%  "

